$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string (row 1)
$ws.Range("A1").Value = "Datos actualizados a 26 de Abril de 2020 a las 03:52"

# --- Corea del Sur (row 34): new case counts ---
$ws.Range("B34").Value = 10728
$ws.Range("C34").Value = 10
$ws.Range("D34").Value = 8717
$ws.Range("E34").Value = 1769
$ws.Range("F34").Value = 55
$ws.Range("G34").Value = 2
$ws.Range("H34").Value = 242

# --- Australia (row 46): new case counts ---
$ws.Range("B46").Value = 6710
$ws.Range("C46").Value = 15
$ws.Range("D46").Value = 5517
$ws.Range("E46").Value = 1112
$ws.Range("F46").Value = 43
$ws.Range("G46").Value = 1
$ws.Range("H46").Value = 81

# --- Jamaica overtakes Tanzania in ranking: swap the two rows (123/124) ---
# Row 123 becomes Jamaica with updated numbers
$ws.Range("A123").Value = "Jamaica"
$ws.Range("B123").Value = 305
$ws.Range("C123").Value = 17
$ws.Range("D123").Value = 28
$ws.Range("E123").Value = 270
$ws.Range("F123").Value = 0
$ws.Range("G123").Value = 0
$ws.Range("H123").Value = 7

# Row 124 becomes Tanzania with its previous (unchanged) numbers
$ws.Range("A124").Value = "Tanzania"
$ws.Range("B124").Value = 299
$ws.Range("C124").Value = 0
$ws.Range("D124").Value = 48
$ws.Range("E124").Value = 241
$ws.Range("F124").Value = 7
$ws.Range("G124").Value = 0
$ws.Range("H124").Value = 10

# --- Paraguay (row 128): new case counts ---
$ws.Range("B128").Value = 228
$ws.Range("C128").Value = 5
$ws.Range("D128").Value = 85
$ws.Range("E128").Value = 134
